# Refactor code structure for improved readability and maintainability
#
# Bumps the "Fecha de ingreso" date column (I4:I6) forward by one day
# and updates the active selection to reflect the single active cell I6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Advance the entry dates in I4, I5 and I6 by one day (46008 -> 46009).
$ws.Range("I4").Value = $ws.Range("I4").Value2() + 1
$ws.Range("I5").Value = $ws.Range("I5").Value2() + 1
$ws.Range("I6").Value = $ws.Range("I6").Value2() + 1

# Update the sheet's active selection from I5:I6 to a single cell I6.
$null = $ws.Range("I6").Select()
